$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the three runs that spell out
#   "4a. T" + "he charge cable will show animation" + " using 2 methods:"
# into a single run "4a. The charge cable will show animation using 2 methods:"
# Find.Execute can match text that spans several runs; replacing that
# matched range with the same (already-concatenated) text collapses it
# into one run.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "4a. T" + "he charge cable will show animation" + " using 2 methods:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "4a. The charge cable will show animation using 2 methods:", 2)
Write-Output ("Change1 Found: " + $found1)

# ---------------------------------------------------------------------
# Change 2: split the single run containing
#   "<U+201C>Screen 2 <U+201C> "
# into three runs:
#   "<U+201C>Screen 2"
#   "<U+201D>"
#   " "
# (the extra inner space is dropped and the second left quote becomes a
# right quote).
# ---------------------------------------------------------------------
$leftQuote  = [char]0x201C
$rightQuote = [char]0x201D

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    $leftQuote + "Screen 2 " + $leftQuote + " ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Change2 Found: " + $found2)

$runStart = $rng2.Start
$runEnd   = $rng2.End

# Current characters (relative to $runStart):
#   0:"  1:S 2:c 3:r 4:e 5:e 6:n 7:sp 8:2 9:sp 10:"  11:sp   (12 chars)
#
# Any edit that changes the *text* of a run can trigger Word to
# re-normalise (merge) adjacent runs that end up with indistinguishable
# formatting, so fix up the textual content FIRST (as one single-run
# edit), and only AFTER that split the now-correct text into separate
# runs by toggling direct formatting on sub-ranges (apply then remove -
# this forces a run boundary without altering the visible text or
# requiring any further text edits that could cause re-merging).

# 1) Normalise the text in one shot: drop the extra inner space and
#    turn the second left-quote into a right-quote.
$whole = $d.Range($runStart, $runEnd)
$whole.Text = $leftQuote + "Screen 2" + $rightQuote + " "

# New characters (relative to $runStart), 11 chars total:
#   0:"  1:S 2:c 3:r 4:e 5:e 6:n 7:sp 8:2 9:"  10:sp

# 2) Split "<U+201C>Screen 2" (0-8) away from the closing quote + space
#    (9-10) tail.
$openingPart = $d.Range($runStart, $runStart + 9)
$openingPart.Bold = 1
$openingPart.Bold = 0

# 3) Split the closing quote (9) away from the trailing space (10).
$closingQuote = $d.Range($runStart + 9, $runStart + 10)
$closingQuote.Bold = 1
$closingQuote.Bold = 0

Write-Output ("Final text: [" + $d.Range($runStart, $runStart + 10).Text + "]")
